# Generate Report for Handback
# Refresh the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps that are stamped each time the
# handback status report is regenerated.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet - Latest HO Xliff Generate Date
$wsOverview.Range("G2").Value = "2016-08-23 17:10:27"

# zh-cn sheet - Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H2").Value = "2016-08-23 17:10:22"
$wsZhCn.Range("K2").Value = "2016-08-23 17:10:41"

# de-de sheet - Correspond Handoff Datetime / Correspond Handback DateTime
$wsDeDe.Range("H2").Value = "2016-08-23 17:10:27"
$wsDeDe.Range("K2").Value = "2016-08-23 17:10:49"
